$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 5): Meal / Mon / Thur / Tues / Weds
$ws.Range("J5").Value = "Meal"
$ws.Range("K5").Value = "Mon"
$ws.Range("L5").Value = "Thur"
$ws.Range("M5").Value = "Tues"
$ws.Range("N5").Value = "Weds"

# Row 6: Breakfast
$ws.Range("J6").Value = "Breakfast"
$ws.Range("K6").Value = "Toast"
$ws.Range("L6").Value = "Toast"
$ws.Range("M6").Value = "Toast"
$ws.Range("N6").Value = "Toast"

# Row 7: Lunch
$ws.Range("J7").Value = "Lunch"
$ws.Range("K7").Value = "Soup"
$ws.Range("L7").Value = "Hotpot"
$ws.Range("M7").Value = "Something Different!"
$ws.Range("N7").Value = "Soup"

# Row 8: Dinner
$ws.Range("J8").Value = "Dinner"
$ws.Range("K8").Value = "Curry"
$ws.Range("L8").Value = "Curry"
$ws.Range("M8").Value = "Curry"
$ws.Range("N8").Value = "Curry"

# Row 9: Midnight Snack
$ws.Range("J9").Value = "Midnight Snack"
$ws.Range("K9").Value = "Shmores"
$ws.Range("L9").Value = "Chocolate"
$ws.Range("M9").Value = "Shmores"
$ws.Range("N9").Value = "Biscuits"

# Copy header style from J5/K5 into the new header cells L5:N5
$ws.Range("K5").Copy()
$ws.Range("L5:N5").PasteSpecial(-4122) # xlPasteFormats

# Remove the old leftover rows (10-16) so the used range shrinks back to I5:N9
$ws.Range("I10:K16").Clear()

